$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Mix-B"
$ws.Range("H1").Value = "Mix-C"
$ws.Range("I1").Value = "Mix-D"

$ws.Range("K15").Select()
